# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the "Type of document" comparison table) gets its
#    table style switched from the deck's custom "Table_0" style to the
#    built-in style {704CDCC9-3988-4697-A960-A90253A62F54}.
#
# 2) The deck's theme ("Integral" / "Red Violet" colour scheme) is swapped
#    for the plain default "Office Theme" colour palette (this is what the
#    underlying theme1.xml content becomes after the edit). Table.Style is
#    read-only in this object model, so table styles must go through
#    Table.ApplyStyle(), and theme colours are changed one swatch at a time
#    through ThemeColorScheme, which is the supported way to repaint the
#    theme's 12-colour palette from script.

$p = $ppt.ActivePresentation

# --- 1. Retarget the slide-5 table's style -------------------------------
$targetStyleId = "{704CDCC9-3988-4697-A960-A90253A62F54}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Recolour the theme from "Red Violet" to the default "Office" -----
function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToVbaRgb($officeColors[$i - 1])
}
